# Regenerate the localization-status report:
#  - the handoff for 340dc998-7864-4936-8b32-2e9540ffa8ee.md completed, producing a
#    new source file (95cb3bac-242a-46e3-9c0e-23f8873970d5.md) with a fresh handoff
#    package (...37ea7e08c1b9fadd3197ff895beb2718c50a9668...) and new handoff
#    timestamps.
#  - the previously-failed handoff row (da68e837-979a-4f2b-9f74-9ba8a4d8b1ab.md /
#    "Handoff transform failed") is gone from this run's report, so its row is
#    removed from every sheet.

$wb = $excel.ActiveWorkbook

$newGuidFile  = "95cb3bac-242a-46e3-9c0e-23f8873970d5.md"
$newGuidBase  = "https://github.com/OpenLocalizationTest/oltest/blob/dffd63f7c14cdcc40dbd42f244dd587a1b050422/e2e/95cb3bac-242a-46e3-9c0e-23f8873970d5.md"
$configAddr   = "https://github.com/OpenLocalizationTest/oltest/blob/dffd63f7c14cdcc40dbd42f244dd587a1b050422/.localization-config"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# The failed-handoff row (old row 3) is dropped; row 4 (.localization-config)
# shifts up into its place.
$ws.Rows.Item(3).Delete()

# New source file name for row 2.
$ws.Range("A2").Value = $newGuidFile

# Rebuild hyperlinks from scratch so ranges/targets line up with the new
# layout (the engine does not auto-shift hyperlink anchors on row delete).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newGuidBase, "", "", $newGuidFile)
$ws.Hyperlinks.Add($ws.Range("A3"), $configAddr, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows.Item(3).Delete()

$ws.Range("A2").Value = $newGuidFile
$ws.Range("C2").Value = "95cb3bac-242a-46e3-9c0e-23f8873970d5.37ea7e08c1b9fadd3197ff895beb2718c50a9668.zh-cn.xlf"
$ws.Range("D2").Value = "2016-01-13 16:06:49"

$zhXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c004468d87b3a31c759a667e1c1a896795bd87f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/95cb3bac-242a-46e3-9c0e-23f8873970d5.37ea7e08c1b9fadd3197ff895beb2718c50a9668.zh-cn.xlf"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newGuidBase, "", "", $newGuidFile)
$ws.Hyperlinks.Add($ws.Range("C2"), $zhXlfAddr, "", "", "95cb3bac-242a-46e3-9c0e-23f8873970d5.37ea7e08c1b9fadd3197ff895beb2718c50a9668.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $configAddr, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows.Item(3).Delete()

$ws.Range("A2").Value = $newGuidFile
$ws.Range("C2").Value = "95cb3bac-242a-46e3-9c0e-23f8873970d5.37ea7e08c1b9fadd3197ff895beb2718c50a9668.de-de.xlf"
$ws.Range("D2").Value = "2016-01-13 16:06:57"

$deXlfAddr = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/606b28171b5e2735e62d324f27e973d363259061/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/95cb3bac-242a-46e3-9c0e-23f8873970d5.37ea7e08c1b9fadd3197ff895beb2718c50a9668.de-de.xlf"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), $newGuidBase, "", "", $newGuidFile)
$ws.Hyperlinks.Add($ws.Range("C2"), $deXlfAddr, "", "", "95cb3bac-242a-46e3-9c0e-23f8873970d5.37ea7e08c1b9fadd3197ff895beb2718c50a9668.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $configAddr, "", "", ".localization-config")
